$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 216
$ws1.Range("F7").Value  = 95
$ws1.Range("F8").Value  = 10201
$ws1.Range("F10").Value = 3516
$ws1.Range("F11").Value = 209
$ws1.Range("F12").Value = 2443
$ws1.Range("F14").Value = 2803
$ws1.Range("F17").Value = 2169
$ws1.Range("F23").Value = 141
$ws1.Range("F26").Value = 221
$ws1.Range("F28").Value = 1317
$ws1.Range("F29").Value = 11
$ws1.Range("F30").Value = 1254
$ws1.Range("F34").Value = 3480
$ws1.Range("F35").Value = 3120
$ws1.Range("F36").Value = 30
$ws1.Range("F38").Value = 1042
$ws1.Range("F39").Value = 400
$ws1.Range("F41").Value = 1292
$ws1.Range("F42").Value = 95
$ws1.Range("F46").Value = 40
$ws1.Range("F47").Value = 8

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 178
$ws2.Range("F7").Value  = 2
$ws2.Range("F8").Value  = 6
$ws2.Range("F16").Value = 177

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 744
$ws3.Range("F5").Value = 1998

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 95
$ws4.Range("F10").Value = 10201
$ws4.Range("F12").Value = 3517
$ws4.Range("F13").Value = 209
$ws4.Range("F17").Value = 2169
$ws4.Range("F21").Value = 141
$ws4.Range("F24").Value = 221
$ws4.Range("F26").Value = 1317
$ws4.Range("F27").Value = 11
$ws4.Range("F28").Value = 1254
$ws4.Range("F30").Value = 2
$ws4.Range("F31").Value = 6
$ws4.Range("F33").Value = 3481
$ws4.Range("F34").Value = 3120
$ws4.Range("F35").Value = 30
$ws4.Range("F36").Value = 1042
$ws4.Range("F39").Value = 400
$ws4.Range("F44").Value = 95
$ws4.Range("F47").Value = 40
$ws4.Range("F48").Value = 8
$ws4.Range("F49").Value = 177
